$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.049002
$ws.Range("H2").Value = 0.147006
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 37.897696
$ws.Range("N2").Value = 113.693088
$ws.Range("O2").Value = 0.5443297838974656
$ws.Range("P2").Value = 0.5443297838974654
$ws.Range("Q2").Value = 1.857062899392
$ws.Range("R2").Value = 16.713566094528
$ws.Range("S2").Value = 0.5443297838974656
$ws.Range("T2").Value = 0.5443297838974654

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.049002
$ws.Range("H3").Value = 0.147006
$ws.Range("O3").Value = 0.1845108538321186
$ws.Range("P3").Value = 0.1845108538321186
$ws.Range("Q3").Value = 0.629486519612
$ws.Range("R3").Value = 5.665378676508
$ws.Range("S3").Value = 0.1845108538321186
$ws.Range("T3").Value = 0.1845108538321186

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.049002
$ws.Range("H4").Value = 0.147006
$ws.Range("M4").Value = 18.829808
$ws.Range("N4").Value = 56.489424
$ws.Range("O4").Value = 0.2704551041696774
$ws.Range("P4").Value = 0.2704551041696774
$ws.Range("Q4").Value = 0.9226982516159999
$ws.Range("R4").Value = 8.304284264544
$ws.Range("S4").Value = 0.2704551041696774
$ws.Range("T4").Value = 0.2704551041696774

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.049002
$ws.Range("H5").Value = 0.147006
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.04903233333333334
$ws.Range("N5").Value = 0.147097
$ws.Range("O5").Value = 0.0007042581007384151
$ws.Range("P5").Value = 0.000704258100738415
$ws.Range("Q5").Value = 0.002402682398
$ws.Range("R5").Value = 0.021624141582
$ws.Range("S5").Value = 0.0007042581007384151
$ws.Range("T5").Value = 0.000704258100738415
